{"js": "// The document has two tables (both styled \"Table\") separated by a spacer\n// paragraph (style \"Body Text\", containing a single space). Re-rendering\n// the Quarto source dropped that stray paragraph, so the second table now\n// directly follows the first one.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet spacer = null;\n\n// Primary approach: the spacer is the paragraph immediately after the\n// first table.\nif (tables.items.length >= 2) {\n  const firstTable = tables.items[0];\n  const candidate = firstTable.getParagraphAfterOrNullObject();\n  candidate.load(\"text,style\");\n  await context.sync();\n\n  if (!candidate.isNullObject && candidate.text.trim() === \"\" && candidate.style === \"Body Text\") {\n    spacer = candidate;\n  }\n}\n\n// Fallback: scan the flat paragraph list (this includes paragraphs that\n// live inside table cells) for an empty \"Body Text\" paragraph, outside of\n// any table, that sits directly between two table-hosted paragraphs.\nif (!spacer) {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const items = paragraphs.items;\n  for (const p of items) {\n    p.load(\"text,style,parentTableOrNullObject\");\n  }\n  await context.sync();\n\n  for (let i = 1; i < items.length - 1; i++) {\n    const para = items[i];\n    const prev = items[i - 1];\n    const next = items[i + 1];\n    const paraInTable = !para.parentTableOrNullObject.isNullObject;\n    const prevInTable = !prev.parentTableOrNullObject.isNullObject;\n    const nextInTable = !next.parentTableOrNullObject.isNullObject;\n\n    if (\n      !paraInTable &&\n      prevInTable &&\n      nextInTable &&\n      para.text.trim() === \"\" &&\n      para.style === \"Body Text\"\n    ) {\n      spacer = para;\n      break;\n    }\n  }\n}\n\nif (spacer) {\n  spacer.delete();\n  await context.sync();\n}\n", "ps1": "# The document has two tables (both styled \"Table\") with a spacer\n# paragraph in between: style \"Body Text\", containing just a single\n# space. Re-rendering the Quarto source dropped that stray paragraph, so\n# the second table now directly follows the first one.\n\n$d = $word.ActiveDocument\n\n$target = $null\n\nif ($d.Tables.Count -ge 2) {\n    # Primary approach: the spacer paragraph is the content strictly\n    # between the end of the first table and the start of the second one.\n    $t1 = $d.Tables.Item(1)\n    $t2 = $d.Tables.Item(2)\n    $gap = $d.Range($t1.Range.End, $t2.Range.Start)\n    if ($gap.Text.Trim() -eq \"\") {\n        $target = $gap\n    }\n}\n\nif ($target -eq $null) {\n    # Fallback: walk the document paragraphs and find the first one that\n    # is outside any table, styled \"Body Text\", and whitespace-only.\n    foreach ($p in $d.Paragraphs) {\n        $isInTable = $p.Range.Information(12)\n        if ((-not $isInTable) -and ($p.Style.NameLocal -eq \"Body Text\") -and ($p.Range.Text.Trim() -eq \"\")) {\n            $target = $p.Range\n            break\n        }\n    }\n}\n\nif ($target -ne $null) {\n    $target.Delete()\n}\n"}
